# Apply updates: add spaces around the hyphen in the "具体时间范围" (E column)
# time-range strings, and bump a couple of "想去人数" (F column) counts by 1,
# across the relevant worksheets (展览, 演出, 全部类型).

$wb = $excel.ActiveWorkbook

function Update-TimeRange {
    param($Worksheet, [string]$CellRef)
    $cell = $Worksheet.Range($CellRef)
    $cell.Value = $cell.Value().Replace("-", " - ")
}

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
Update-TimeRange $wsExhibit "E2"
Update-TimeRange $wsExhibit "E3"
$wsExhibit.Range("F3").Value = 2122
Update-TimeRange $wsExhibit "E4"
Update-TimeRange $wsExhibit "E5"
$wsExhibit.Range("F5").Value = 1350
Update-TimeRange $wsExhibit "E6"

# --- Sheet "演出" ---
$wsShow = $wb.Worksheets.Item("演出")
Update-TimeRange $wsShow "E2"
Update-TimeRange $wsShow "E3"

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
Update-TimeRange $wsAll "E2"
Update-TimeRange $wsAll "E3"
$wsAll.Range("F3").Value = 2122
Update-TimeRange $wsAll "E4"
Update-TimeRange $wsAll "E5"
Update-TimeRange $wsAll "E6"
Update-TimeRange $wsAll "E7"
$wsAll.Range("F7").Value = 1350
Update-TimeRange $wsAll "E8"
